$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "Status" column before the old "Jan_2026" column (D) ---
# This shifts old D..H (Jan_2026..QoQ) to E..I, matching the diff's column layout.
$ws.Columns("D").Insert()
$ws.Range("D1").Value = "Status"

# Old F1 ("Nov_2025") shifted to G1; the new header set renames it to "Oct_2025".
$ws.Range("G1").Value = "Oct_2025"

# --- Insert 6 new rows at the bottom (21..26) for the newly added holdings ---
$ws.Rows("21:26").Insert()

# --- Full target data for rows 2..26: ISIN, Stock Name, Status, Jan_2026, Dec_2025, Oct_2025, MoM, QoQ ---
$data = @(
    @("INE406A01037", "Aurobindo Pharma Limited", "Adding Consistently", 10.101106, 9.993342, 8.983580999999999, 0.1077639999999995, 1.117525000000001),
    @("INE775A01035", "Samvardhana Motherson International Ltd", "Reducing", 9.837403, 10.041152, 9.186685000000001, -0.2037490000000002, 0.6507179999999995),
    @("INE423A01024", "Adani Enterprises Limited", "Reducing", 9.614039, 9.733378, 0, -0.1193390000000001, 9.614039),
    @("INE768C01028", "Zydus Wellness Ltd", "Adding Consistently", 8.306566, 7.692854, 7.580929, 0.6137120000000005, 0.7256369999999999),
    @("INE045A01017", "Ador Welding Limited", "Adding Consistently", 6.21963, 5.795094, 5.174521, 0.4245360000000007, 1.045109),
    @("INE917I01010", "Bajaj Auto Limited", "Adding Consistently", 6.031482, 5.362379, 4.528778, 0.6691029999999998, 1.502704),
    @("INE364U01010", "Adani Green Energy Limited", "Reducing Consistently", 5.856657, 6.370469, 7.362987, -0.5138119999999997, -1.50633),
    @("INE206N01018", "Ravindra Energy Limited", "Adding Consistently", 5.600932, 5.541504, 5.096678, 0.05942800000000048, 0.5042540000000004),
    @("INE942C01045", "Gujarat Themis Biosyn Ltd", "Reducing Consistently", 5.053164, 6.308265, 5.894732, -1.255101, -0.8415680000000005),
    @("INE180C01042", "Capri Global Capital Limited", "Adding Consistently", 4.176886, 3.968793, 3.910511, 0.2080929999999999, 0.2663749999999996),
    @("INE931S01010", "Adani Energy Solutions Limited", "Reducing", 3.793646, 3.977906, 3.388341, -0.1842600000000001, 0.4053049999999998),
    @("INE331A01037", "The Ramco Cements Limited", "Fresh Entry", 3.159637, 0, 0, 3.159637, 3.159637),
    @("INE034A01011", "Arvind Limited", "Complete Exit", 0, 0, 2.82423, 0, -2.82423),
    @("INE290A01027", "Nahar Spinning Mills Limited", "Complete Exit", 0, 0.577651, 0.545767, -0.577651, -0.545767),
    @("INE373A01013", "BASF India Ltd", "Complete Exit", 0, 0, 3.053269, 0, -3.053269),
    @("INE258G01013", "Sumitomo Chemical India Limited", "Complete Exit", 0, 0, 2.079425, 0, -2.079425),
    @("INE182A01018", "Pfizer Ltd", "Complete Exit", 0, 0, 2.618242, 0, -2.618242),
    @("INE548A01028", "HFCL Limited", "Complete Exit", 0, 0, 2.2633, 0, -2.2633),
    @("INE669C01036", "Tech Mahindra Limited", "Complete Exit", 0, 1.336715, 0, -1.336715, 0),
    @("INE171Z01026", "Bharat Dynamics Limited", "Complete Exit", 0, 5.330176, 0, -5.330176, 0),
    @("INE769A01020", "Aarti Industries Ltd", "Complete Exit", 0, 0, 2.600919, 0, -2.600919),
    @("INE14LE01019", "Aditya Birla Lifestyle Brands Limited", "Complete Exit", 0, 0, 0.850031, 0, -0.850031),
    @("INE0CZ201020", "ANTHEM BIOSCIENCES LIMITED", "Complete Exit", 0, 0, 5.026891, 0, -5.026891),
    @("INE0BS701011", "Premier Energies Limited", "Complete Exit", 0, 0, 6.02133, 0, -6.02133),
    @("INE019C01026", "Himadri Speciality Chemical Limited", "Complete Exit", 0, 0, 2.878716, 0, -2.878716)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = "quant Manufacturing Fund"
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
    $ws.Cells.Item($r, 9).Value = $row[7]
    $r = $r + 1
}
